$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): keep A1/C1 text as-is, retype the rest ---
# (order matches the shared-string insertion order of the target file)
$ws.Range("E1").Value = "kind"
$ws.Range("D1").Value = "id"
$ws.Range("A2").Value = "jorge"
$ws.Range("B1").Value = "locacalizacion"

# --- Row 2 (values) ---
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1

# --- New font for B2 (monospace, olive-green, size 10) ---
$ws.Range("B2").Font.Name = "Courier New"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").Font.Color = 5867370
$ws.Range("B2").VerticalAlignment = -4108

# --- Drop the now-unused trailing columns F:I ---
$ws.Range("F1:I2").ClearContents()

# --- Selection / view matches the shrunk used range ---
$ws.Range("A1:E2").Select()
